$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the B2:D9 range of numeric values to 0, matching the diff (also fills
# the previously-empty C2 cell with a numeric 0).
$ws.Range("B2:D9").Value = 0
